# ---------------------------------------------------------------------------
# Applies the "Working file with executing using pom.xml" commit:
#   1. Inserts a brand new worksheet ("Sheet4") between Sheet2 and Sheet3,
#      holding a small test-run summary / defects table.
#   2. Appends more rows of test-case data to Sheet2 (TC6..TC9 reworked).
#   3. Edits Sheet3's existing rows (TC3..TC5 reshuffled) and removes the
#      last row (old TC6), keeping Sheet3 as the active/selected tab.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new worksheet right before Sheet3 (Sheet3 is the sheet that
#    was active/selected in the original workbook) and populate it with the
#    defects / test-run summary table.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$newSheet = $wb.Worksheets.Add($ws3)

$newSheet.Range("A1").Value = "Total Number of Test cases "
$newSheet.Range("B1").Value = "Total number of Automated Scripts "
$newSheet.Range("C1").Value = "Number of test cases failed "
$newSheet.Range("D1").Value = "Number of test cases passed"
$newSheet.Range("E1").Value = "Defects "

$newSheet.Range("A2").Value = 17
$newSheet.Range("B2").Value = 10
$newSheet.Range("C2").Value = 3
$newSheet.Range("E2").Value = "In firstname and last name taking single characters "

$newSheet.Range("E3").Value = "first name and last name taking single numbers "
$newSheet.Range("E4").Value = "password field is taking more number of characters"
$newSheet.Range("E5").Value = "If we numbers and characters in the email field some times it was not giving error"
$newSheet.Range("E6").Value = "phone number is accepting single number . "

$newSheet.Columns.Item(1).ColumnWidth = 24
$newSheet.Columns.Item(2).ColumnWidth = 30.90625
$newSheet.Columns.Item(3).ColumnWidth = 24.08984375
$newSheet.Columns.Item(4).ColumnWidth = 24.81640625
$newSheet.Columns.Item(5).ColumnWidth = 53.6328125
$newSheet.Columns.Item(6).ColumnWidth = 22.26953125

$newSheet.Application.ActiveWindow.ScrollColumn = 4
$newSheet.Range("F1").Select()

# Re-select Sheet3 (the sheet that was active before the insert) so it keeps
# the "tabSelected" / activeTab state, matching the original file's intent.
$ws3Again = $wb.Worksheets.Item("Sheet3")
$ws3Again.Select()

# ---------------------------------------------------------------------------
# 2) Sheet2: rework the TC6 row, add new TC7/TC8/TC9 rows.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 7 (TC6) - overwrite with the "verify all the fields with numbers" data.
$ws2.Range("B7").Value = "Verify all the fields with numbers"
$ws2.Range("D7").Value = 1
$ws2.Range("E7").Value = 2
$ws2.Range("F7").Value = 3
$ws2.Range("G7").Value = 4
$ws2.Range("H7").Value = 5
$ws2.Range("I7").Value = 5

# Row 8 (TC7) - new content, clear the now-unused trailing cells.
$ws2.Range("B8").Value = "verify all the fields with numbers "
$ws2.Range("D8").Value = 2343
$ws2.Range("E8").Value = 34353
$ws2.Range("F8").ClearContents()
$ws2.Range("G8").ClearContents()
$ws2.Range("H8").ClearContents()
$ws2.Range("I8").ClearContents()

# Row 9 (TC8) - new content.
$ws2.Range("B9").Value = "verify email field with numbers and chracters"
$ws2.Range("D9").ClearContents()
$ws2.Range("E9").Value = "123asdertfffm34557"
$ws2.Range("F9").ClearContents()
$ws2.Range("G9").ClearContents()
$ws2.Range("H9").ClearContents()
$ws2.Range("I9").ClearContents()

# Row 10 (TC9) - brand new row.
$ws2.Range("A10").Value = "TC9"
$ws2.Range("B10").Value = "Verify password fields with different passwords"
$ws2.Range("C10").Value = "Yes"
$ws2.Range("D10").Value = "abc"
$ws2.Range("E10").Value = "def"
$ws2.Range("F10").Value = "def@gmail.com"
$ws2.Hyperlinks.Add($ws2.Range("F10"), "mailto:def@gmail.com")
$ws2.Range("G10").Value = 3345678901
$ws2.Range("H10").Value = "jkl"
$ws2.Range("I10").Value = "lkj"

$ws2.Columns.Item(4).ColumnWidth = 15.1796875

$ws2.Range("A11").Select()

# ---------------------------------------------------------------------------
# 3) Sheet3 (the original sheet, still named Sheet3, now the last tab):
#    reshuffle the TC3/TC4/TC5 rows and drop the old TC6 row entirely.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")

# Row 4 (TC3) now carries what used to be the TC4 data.
$ws3.Range("D4").Value = "Submit Request with invoice Request tupe"
$ws3.Range("G4").Value = "Invoice Request"
$ws3.Range("H4").Value = 4568
$ws3.Range("I4").Value = "I want my Invoice"

# Row 5 (TC4) now carries what used to be the TC5 data, clear the rest.
$ws3.Range("D5").Value = "Verify with null data"
$ws3.Range("E5").Value = "no"
$ws3.Range("F5").ClearContents()
$ws3.Range("G5").ClearContents()
$ws3.Range("H5").ClearContents()
$ws3.Range("I5").ClearContents()

# Row 6 (TC5) loses its TestCase text - only the login/TCID columns remain.
$ws3.Range("D6").ClearContents()

# Drop the old TC6 row (row 7) completely.
$ws3.Range("A7").EntireRow.Delete()

# The engine's "delete row" does not prune hyperlinks bound to the deleted
# row, so rebuild the hyperlink list for the rows that remain.
$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "mailto:srikanthtesting100@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "mailto:srikanthtesting100@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "mailto:srikanthtesting100@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("A5"), "mailto:srikanthtesting100@gmail.com")
$ws3.Hyperlinks.Add($ws3.Range("A6"), "mailto:srikanthtesting100@gmail.com")

$ws3.Range("B10").Select()
